$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.239.63"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.81%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'2.630.28"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.08%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  -0.11%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'597.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.24%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'152.23"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +1.34%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.07%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.590"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.28%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.115"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +5.35%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'5.84"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +2.65%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.396"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +3.41%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("E12").Value = "'  +1.14%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'28.08"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.18%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'3.099.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.00%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.0000173"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +15.33%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'63.980.34"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.70%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'2.595.90"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.28%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'12.32"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.11%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'4.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +2.90%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'350.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.72%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'7.13"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +3.64%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("E22").Value = "'  +0.15%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'67.69"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +1.94%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'1.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -3.07%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'9.26"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.71%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'1.68"
$ws.Range("D26").Style = "Normal"

$ws.Range("D27").Value = "'8.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.51%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'550.07"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -3.76%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'0.163"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.69%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'0.998"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.24%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'0.0₃0911"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +7.88%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'2.07"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +1.28%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'1.84"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +5.27%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'5.49"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +4.24%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'6.19"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +1.27%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").Value = "'  +3.70%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'165.35"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -2.16%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'20.20"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +4.22%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'2.01"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +2.72%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.997"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.25%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("E41").Value = "'  -0.01%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'169.40"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.42%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'4.12"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +4.76%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'23.34"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +9.04%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("B45").Value = "'dogwifhat"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'2.22"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +11.81%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("B46").Value = "'Hedera"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'0.0588"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -2.26%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.640"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.89%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'0.0254"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +1.66%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'0.0974"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.85%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'19.36"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.77%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("E51").Value = "'  +19.40%  "
$ws.Range("E51").Style = "Normal"
